$d = $word.ActiveDocument

# Helper: returns a fresh Range covering the paragraph that contains the
# given (unique) anchor text. Re-locating it each time keeps the range
# valid even after earlier edits shift character offsets.
function Get-ParaRange($anchorText) {
    $r = $d.Content
    $r.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $p = $r.Paragraphs(1).Range
    return $d.Range($p.Start, $p.End)
}

$anchor = "As a tournament director, I would like to be able to create a tournament by scheduling it,"

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that currently sits right
#    after "... after searching for them" in the anonymous-user /
#    judge's-philosophy paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Normalize the tournament-director sentence tail into a single
#    run so we can rebuild the run/highlight structure cleanly.
#    (All Find operations below are scoped to this one paragraph so
#    they can't accidentally match similar text elsewhere in the doc.)
# ------------------------------------------------------------------
$full = " uploading logistics for the tournament in .pdf or .docx format, and allow teams to be entered"
$p = Get-ParaRange($anchor)
$p.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, $full, 2) | Out-Null

# ------------------------------------------------------------------
# 3) Highlight everything after the leading space in green.
# ------------------------------------------------------------------
$p = Get-ParaRange($anchor)
$p.Find.Execute("uploading logistics for the tournament in .pdf or .docx format, and allow teams to be entered", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p.Font.HighlightColorIndex = 4

# ------------------------------------------------------------------
# 4) Re-introduce the run break between "...tournament in " and
#    ".pdf or ..." (matching the original document's run boundary)
#    by toggling the highlight color off/on across that sub-range.
# ------------------------------------------------------------------
$p = Get-ParaRange($anchor)
$p.Find.Execute("uploading logistics for the tournament in ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p.Font.HighlightColorIndex = 7
$p.Font.HighlightColorIndex = 4

# ------------------------------------------------------------------
# 5) Re-introduce the run break between ".pdf or " and
#    ".docx format, and allow teams to be entered".
# ------------------------------------------------------------------
$p = Get-ParaRange($anchor)
$p.Find.Execute(".pdf or ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p.Font.HighlightColorIndex = 7
$p.Font.HighlightColorIndex = 4

# ------------------------------------------------------------------
# 6) Re-introduce the run break between ".docx format, and " and
#    "allow teams to be entered".
# ------------------------------------------------------------------
$p = Get-ParaRange($anchor)
$p.Find.Execute(".docx format, and ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p.Font.HighlightColorIndex = 7
$p.Font.HighlightColorIndex = 4

# ------------------------------------------------------------------
# 7) Insert the "_GoBack" bookmark (collapsed) right before
#    "allow teams to be entered".
# ------------------------------------------------------------------
$p = Get-ParaRange($anchor)
$p.Find.Execute("allow teams to be entered", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $d.Range($p.Start, $p.Start)
$d.Bookmarks.Add("_GoBack", $ins)
